$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.278.37'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '2.363.45'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.13'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.94'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  +5.54%  '
$ws.Range("E9").Value = '  +4.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.53'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.21%  '
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.06'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").Value = '2.781.62'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").Value = '58.197.17'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("E16").Value = '  +2.65%  '
$ws.Range("D17").Value = '2.353.40'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.91'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.03%  '
$ws.Range("E19").Value = '  +2.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '331.46'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.86'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.00%  '
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.60'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.47%  '
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.27'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.33'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -5.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.91'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.39'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("E33").Value = '  -2.90%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.20'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.405'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '141.91'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.76%  '
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '288.77'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0952'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.84%  '
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("E46").Value = '  +1.05%  '
$ws.Range("E47").Value = '  +2.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.386'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.54'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.14%  '
